$d = $word.ActiveDocument

# --- Text replacements (Find to locate, then set Range.Text to avoid smart-quote autocorrect) ---
$r0 = $d.Content
$found0 = $r0.Find.Execute("Unveiling the Enigma of Dark Matter")
if (-not $found0) { Write-Host "FAILED find 0: Unveiling the Enigma of Dark Matter" } else { $r0.Text = "Chemistry: The Realm of Elements and Interactions" }
$r1 = $d.Content
$found1 = $r1.Find.Execute("Dr. Amelia Grayson")
if (-not $found1) { Write-Host "FAILED find 1: Dr. Amelia Grayson" } else { $r1.Text = "Erica Peterson" }
$r2 = $d.Content
$found2 = $r2.Find.Execute("amelia")
if (-not $found2) { Write-Host "FAILED find 2: amelia" } else { $r2.Text = "erica" }
$r3 = $d.Content
$found3 = $r3.Find.Execute("grayson@stellarobservatory")
if (-not $found3) { Write-Host "FAILED find 3: grayson@stellarobservatory" } else { $r3.Text = "p.chemistry@eduworld" }
$r4 = $d.Content
$found4 = $r4.Find.Execute("In the cosmic tapestry of the universe, dark matter remains an enigmatic entity, an invisible force shaping the galaxies' structure and dynamics")
if (-not $found4) { Write-Host "FAILED find 4: In the cosmic tapestry of the universe, dark matter remains " } else { $r4.Text = "In the vast panorama of scientific inquiry, chemistry stands as a captivating discipline that seeks to decipher the intricate world of elements and their interactions" }
$r5 = $d.Content
$found5 = $r5.Find.Execute(" Its gravitational influence orchestrates the motion of stars, while its shadowy nature eludes direct observation")
if (-not $found5) { Write-Host "FAILED find 5:  Its gravitational influence orchestrates the motion of star" } else { $r5.Text = " Chemistry unravels the secrets of matter, exploring the fundamental building blocks of the universe and the forces that govern their behavior" }
$r6 = $d.Content
$found6 = $r6.Find.Execute(" Unraveling the secrets of dark matter is not merely an academic pursuit but a fundamental endeavor that holds the key to understanding the very fabric of space and time. It's a quest that has captivated the scientific community, requiring innovative approaches and interdisciplinary collaborations")
if (-not $found6) { Write-Host "FAILED find 6:  Unraveling the secrets of dark matter is not merely an acad" } else { $r6.Text = " It unveils the enchanting tapestry of chemical reactions, where atoms dance in a mesmerizing choreography, forging new substances with remarkable properties" }
$r7 = $d.Content
$found7 = $r7.Find.Execute("Dark matter makes up over 85% of the universe's mass, yet its true identity remains veiled")
if (-not $found7) { Write-Host "FAILED find 7: Dark matter makes up over 85% of the universe's mass, yet it" } else { $r7.Text = "From the colossal stars that illuminate the night sky to the microscopic organisms that inhabit the depths of the oceans, chemistry underpins the very essence of life" }
$r8 = $d.Content
$found8 = $r8.Find.Execute(" The most prevalent hypothesis posits the existence of weakly interacting massive particles (WIMPs), elusive entities that evade detection due to their feeble interactions with ordinary matter")
if (-not $found8) { Write-Host "FAILED find 8:  The most prevalent hypothesis posits the existence of weakl" } else { $r8.Text = " It governs the intricate symphony of biochemical processes that occur within living cells, orchestrating the delicate balance that sustains life" }
$r9 = $d.Content
$found9 = $r9.Find.Execute(" The quest for WIMPs has encompassed a wide array of experiments, from underground laboratories to satellite-based missions. Despite these efforts, the elusive particles continue to elude our grasp, prompting scientists to explore alternative candidates such as axions and sterile neutrinos")
if (-not $found9) { Write-Host "FAILED find 9:  The quest for WIMPs has encompassed a wide array of experim" } else { $r9.Text = " Chemistry is omnipresent, its influence discernible in the air we breathe, the food we consume, and the myriad materials that shape our world" }
$r10 = $d.Content
$found10 = $r10.Find.Execute("The study of dark matter transcends the realm of mere scientific curiosity")
if (-not $found10) { Write-Host "FAILED find 10: The study of dark matter transcends the realm of mere scient" } else { $r10.Text = "As we delve into the realm of chemistry, we embark on an exhilarating journey of discovery" }
$r11 = $d.Content
$found11 = $r11.Find.Execute(" Its implications ripple across various fields, promising breakthroughs in astrophysics, cosmology, and fundamental physics")
if (-not $found11) { Write-Host "FAILED find 11:  Its implications ripple across various fields, promising br" } else { $r11.Text = " We uncover the fundamental principles that govern the behavior of atoms and molecules, unraveling the secrets of chemical reactions and the properties of substances" }
$r12 = $d.Content
$found12 = $r12.Find.Execute(" Solving the dark matter puzzle could revolutionize our understanding of the universe's evolution, the nature of gravity, and the ultimate fate of our cosmos. With each passing day, scientists inch closer to unveiling the enigma of dark matter, bringing us tantalizingly close to unraveling one of the greatest mysteries in the universe")
if (-not $found12) { Write-Host "FAILED find 12:  Solving the dark matter puzzle could revolutionize our unde" } else { $r12.Text = " Through hands-on experimentation and theoretical exploration, we gain insights into the profound impact of chemistry on our lives and the world around us" }
$r13 = $d.Content
$found13 = $r13.Find.Execute("Dark matter, an enigmatic entity comprising over 85% of the universe's mass, continues to perplex scientists with its elusive nature")
if (-not $found13) { Write-Host "FAILED find 13: Dark matter, an enigmatic entity comprising over 85% of the " } else { $r13.Text = "Chemistry, a captivating branch of science, unlocks the mysteries of matter, elements, and their interactions" }
$r14 = $d.Content
$found14 = $r14.Find.Execute(" Despite extensive research and experimental efforts, its true identity remains shrouded in mystery")
if (-not $found14) { Write-Host "FAILED find 14:  Despite extensive research and experimental efforts, its tr" } else { $r14.Text = " It unveils the underlying principles that govern the behavior of atoms and molecules, shedding light on the intricate symphony of chemical reactions and the extraordinary properties of substances" }
$r15 = $d.Content
$found15 = $r15.Find.Execute(" Leading hypotheses propose weakly interacting massive particles (WIMPs), axions, and sterile neutrinos as potential candidates. Unraveling the dark matter enigma holds profound implications for astrophysics, cosmology, and fundamental physics, promising breakthroughs in our understanding of the universe's evolution, the nature of gravity, and the ultimate fate of our cosmos")
if (-not $found15) { Write-Host "FAILED find 15:  Leading hypotheses propose weakly interacting massive parti" } else { $r15.Text = " Chemistry's influence extends far beyond the laboratory, as it plays a pivotal role in life processes, industry, and the materials that shape our world" }

# --- Append a new empty paragraph at the end of the document body ---
$endPos = $d.Content.End
$endRng = $d.Range($endPos, $endPos)
$endRng.InsertParagraphAfter()

Write-Host "DONE"